$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.593.75'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.698.10'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.64'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.80%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4038'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.73%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.547'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +7.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '54.80'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +13.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.001'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08803'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.278'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +11.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.34'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001331'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.658'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +6.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.709.39'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '101.16'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07039'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.61'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.899'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.11'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.591.55'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.957'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +7.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.342'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.38'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.61'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.238'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.75'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.675'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +28.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.114'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.889.37'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.512'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +16.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08558'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.81%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.21'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +8.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.986'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2747'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +4.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.77'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02795'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +11.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09039'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.471'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7766'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7287'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.59'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +5.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.508'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.193'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.85%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.56'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.300'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +15.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.08009'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.30%  '
